$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-9 (A:T)
$data = @(
  @("ECs","Hc","C5ar2","ECs",1,0.3333333333333333,0.4782376666666666,1.434713,0.6428035065113115,0.6428035065113116,2,0.6666666666666666,24.12906,72.38718,0.7578756968848996,0.7578756968848996,11.53942535326,103.85482817934,0.4871651554573173,0.4871651554573174),
  @("ECs","Hc","C5ar2","FAPs",1,0.3333333333333333,0.4782376666666666,1.434713,0.6428035065113115,0.6428035065113116,2,0.6666666666666666,0.6016783333333333,1.805035,0.01889826566702329,0.01889826566702329,0.2877452422172222,2.589707179955,0.0121478714377449,0.0121478714377449),
  @("ECs","Hc","C5ar2","M2",1,0.3333333333333333,0.4782376666666666,1.434713,0.6428035065113115,0.6428035065113116,3,1,7.067793666666666,21.203381,0.2219941038135625,0.2219941038135625,3.38008515162811,30.42076636465299,0.1426985883561941,0.1426985883561941),
  @("ECs","Hc","C5ar2","sCs",1,0.3333333333333333,0.4782376666666666,1.434713,0.6428035065113115,0.6428035065113116,2,0.6666666666666666,0.039222,0.117666,0.001231933634514545,0.001231933634514545,0.018757437762,0.168816939858,0.0007918912600551741,0.0007918912600551742),
  @("FAPs","Hc","C5ar2","ECs",3,1,0.2657496666666667,0.797249,0.3571964934886884,0.3571964934886885,2,0.6666666666666666,24.12906,72.38718,0.7578756968848996,0.7578756968848996,6.412289651979999,57.71060686782,0.2707105414275823,0.2707105414275823),
  @("FAPs","Hc","C5ar2","FAPs",3,1,0.2657496666666667,0.797249,0.3571964934886884,0.3571964934886885,2,0.6666666666666666,0.6016783333333333,1.805035,0.01889826566702329,0.01889826566702329,0.1598958165238889,1.439062348715,0.006750394229278388,0.006750394229278389),
  @("FAPs","Hc","C5ar2","M2",3,1,0.2657496666666667,0.797249,0.3571964934886884,0.3571964934886885,3,1,7.067793666666666,21.203381,0.2219941038135625,0.2219941038135625,1.878263810985444,16.904374298869,0.07929551545736842,0.07929551545736842),
  @("FAPs","Hc","C5ar2","sCs",3,1,0.2657496666666667,0.797249,0.3571964934886884,0.3571964934886885,2,0.6666666666666666,0.039222,0.117666,0.001231933634514545,0.001231933634514545,0.010423233426,0.09380910083399999,0.000440042374459371,0.0004400423744593711)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $rowData = $data[$i]
  for ($col = 1; $col -le $rowData.Count; $col++) {
    $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
  }
}
